{"js": "// Replace the 100 \"A\u00d7B=C\" answer strings in the 20x5 practice table with\n// their updated values, in row-major (reading) order, matching the diff.\nconst newValues = [\"89\u00d799=8811\", \"63\u00d763=3969\", \"34\u00d732=1088\", \"85\u00d779=6715\", \"78\u00d740=3120\", \"27\u00d733=891\", \"84\u00d715=1260\", \"16\u00d762=992\", \"13\u00d739=507\", \"93\u00d725=2325\", \"52\u00d791=4732\", \"74\u00d772=5328\", \"12\u00d727=324\", \"46\u00d760=2760\", \"17\u00d744=748\", \"43\u00d737=1591\", \"100\u00d719=1900\", \"83\u00d778=6474\", \"44\u00d7100=4400\", \"77\u00d7100=7700\", \"17\u00d772=1224\", \"92\u00d795=8740\", \"89\u00d766=5874\", \"82\u00d748=3936\", \"45\u00d742=1890\", \"22\u00d746=1012\", \"73\u00d724=1752\", \"44\u00d752=2288\", \"93\u00d794=8742\", \"67\u00d715=1005\", \"60\u00d743=2580\", \"68\u00d746=3128\", \"19\u00d798=1862\", \"65\u00d767=4355\", \"26\u00d723=598\", \"12\u00d713=156\", \"96\u00d784=8064\", \"25\u00d719=475\", \"64\u00d717=1088\", \"20\u00d727=540\", \"67\u00d778=5226\", \"93\u00d743=3999\", \"60\u00d770=4200\", \"20\u00d731=620\", \"36\u00d779=2844\", \"20\u00d724=480\", \"86\u00d798=8428\", \"36\u00d732=1152\", \"84\u00d715=1260\", \"43\u00d749=2107\", \"44\u00d765=2860\", \"21\u00d770=1470\", \"79\u00d748=3792\", \"92\u00d784=7728\", \"36\u00d725=900\", \"16\u00d735=560\", \"18\u00d760=1080\", \"52\u00d776=3952\", \"22\u00d797=2134\", \"65\u00d779=5135\", \"47\u00d781=3807\", \"23\u00d795=2185\", \"56\u00d723=1288\", \"40\u00d714=560\", \"86\u00d714=1204\", \"44\u00d729=1276\", \"69\u00d783=5727\", \"23\u00d712=276\", \"34\u00d760=2040\", \"30\u00d729=870\", \"53\u00d711=583\", \"75\u00d711=825\", \"89\u00d771=6319\", \"21\u00d721=441\", \"34\u00d779=2686\", \"98\u00d729=2842\", \"77\u00d720=1540\", \"91\u00d730=2730\", \"36\u00d744=1584\", \"43\u00d782=3526\", \"64\u00d728=1792\", \"50\u00d743=2150\", \"58\u00d7100=5800\", \"77\u00d761=4697\", \"23\u00d764=1472\", \"35\u00d752=1820\", \"33\u00d723=759\", \"18\u00d758=1044\", \"52\u00d755=2860\", \"59\u00d763=3717\", \"26\u00d756=1456\", \"28\u00d774=2072\", \"48\u00d746=2208\", \"79\u00d762=4898\", \"79\u00d786=6794\", \"42\u00d740=1680\", \"10\u00d765=650\", \"78\u00d725=1950\", \"47\u00d718=846\", \"80\u00d784=6720\"];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.rows.load(\"items\");\nawait context.sync();\n\nconst rows = table.rows.items;\nfor (let r = 0; r < rows.length; r++) {\n  rows[r].cells.load(\"items\");\n}\nawait context.sync();\n\n// Walk the table in row-major (reading) order, writing each updated\n// equation string into the corresponding cell.\nlet idx = 0;\nfor (let r = 0; r < rows.length && idx < newValues.length; r++) {\n  const cells = rows[r].cells.items;\n  for (let c = 0; c < cells.length && idx < newValues.length; c++) {\n    cells[c].value = newValues[idx];\n    idx++;\n  }\n}\nawait context.sync();\n", "ps1": "# Replace the 100 \"A\u00d7B=C\" answer strings in the 20x5 practice table with\n# their updated values, in row-major (reading) order, matching the diff.\n$newValues = @(\"89\u00d799=8811\", \"63\u00d763=3969\", \"34\u00d732=1088\", \"85\u00d779=6715\", \"78\u00d740=3120\", \"27\u00d733=891\", \"84\u00d715=1260\", \"16\u00d762=992\", \"13\u00d739=507\", \"93\u00d725=2325\", \"52\u00d791=4732\", \"74\u00d772=5328\", \"12\u00d727=324\", \"46\u00d760=2760\", \"17\u00d744=748\", \"43\u00d737=1591\", \"100\u00d719=1900\", \"83\u00d778=6474\", \"44\u00d7100=4400\", \"77\u00d7100=7700\", \"17\u00d772=1224\", \"92\u00d795=8740\", \"89\u00d766=5874\", \"82\u00d748=3936\", \"45\u00d742=1890\", \"22\u00d746=1012\", \"73\u00d724=1752\", \"44\u00d752=2288\", \"93\u00d794=8742\", \"67\u00d715=1005\", \"60\u00d743=2580\", \"68\u00d746=3128\", \"19\u00d798=1862\", \"65\u00d767=4355\", \"26\u00d723=598\", \"12\u00d713=156\", \"96\u00d784=8064\", \"25\u00d719=475\", \"64\u00d717=1088\", \"20\u00d727=540\", \"67\u00d778=5226\", \"93\u00d743=3999\", \"60\u00d770=4200\", \"20\u00d731=620\", \"36\u00d779=2844\", \"20\u00d724=480\", \"86\u00d798=8428\", \"36\u00d732=1152\", \"84\u00d715=1260\", \"43\u00d749=2107\", \"44\u00d765=2860\", \"21\u00d770=1470\", \"79\u00d748=3792\", \"92\u00d784=7728\", \"36\u00d725=900\", \"16\u00d735=560\", \"18\u00d760=1080\", \"52\u00d776=3952\", \"22\u00d797=2134\", \"65\u00d779=5135\", \"47\u00d781=3807\", \"23\u00d795=2185\", \"56\u00d723=1288\", \"40\u00d714=560\", \"86\u00d714=1204\", \"44\u00d729=1276\", \"69\u00d783=5727\", \"23\u00d712=276\", \"34\u00d760=2040\", \"30\u00d729=870\", \"53\u00d711=583\", \"75\u00d711=825\", \"89\u00d771=6319\", \"21\u00d721=441\", \"34\u00d779=2686\", \"98\u00d729=2842\", \"77\u00d720=1540\", \"91\u00d730=2730\", \"36\u00d744=1584\", \"43\u00d782=3526\", \"64\u00d728=1792\", \"50\u00d743=2150\", \"58\u00d7100=5800\", \"77\u00d761=4697\", \"23\u00d764=1472\", \"35\u00d752=1820\", \"33\u00d723=759\", \"18\u00d758=1044\", \"52\u00d755=2860\", \"59\u00d763=3717\", \"26\u00d756=1456\", \"28\u00d774=2072\", \"48\u00d746=2208\", \"79\u00d762=4898\", \"79\u00d786=6794\", \"42\u00d740=1680\", \"10\u00d765=650\", \"78\u00d725=1950\", \"47\u00d718=846\", \"80\u00d784=6720\")\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$rows = $t.Rows.Count\n$cols = $t.Columns.Count\n\n# Walk the table in row-major (reading) order, writing each updated\n# equation string into the corresponding cell's range.\n$idx = 0\nfor ($r = 1; $r -le $rows; $r++) {\n  for ($c = 1; $c -le $cols; $c++) {\n    if ($idx -ge $newValues.Length) { break }\n    $cell = $t.Cell($r, $c)\n    $cell.Range.Text = $newValues[$idx]\n    $idx++\n  }\n}\n"}
